$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell for the new "MatchingInterests" column
$ws.Range("S1").Value = "MatchingInterests"
$ws.Range("R1").Copy()
$ws.Range("S1").PasteSpecial(-4122)

# Computed matching-interests values for each data row
$ws.Range("S2").Value = "Handicrafts,History,Technology,Physics,"
$ws.Range("S3").Value = "Philosophy,English,Social Studies,"
$ws.Range("S4").Value = "Gardening,Physics,"
$ws.Range("S5").Value = "Gardening,Emotional Regulation,Technology,Physics,"
$ws.Range("S6").Value = "Cooking,Baking,Apps,Chemistry,"
